# Generate Report for handback
#
# This script brings the "before" localization-status workbook up to date
# with a fresh handback report:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" on the Overview sheet and on every
#     per-locale sheet.
#   - Each per-locale sheet (zh-cn, de-de) grows two new columns worth of
#     data per data row: "Latest Target File" (E) and "Latest Handback
#     File" (F), each populated + hyperlinked like their sibling columns.
#   - The "Latest Handback DateTime" column (G) moves from the placeholder
#     "0001-01-01 00:00:00" to a real handback timestamp.
#   - The bottom ".localization-config" row's placeholder datetime/reason
#     cells are (re)written to keep everything consistent.

$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

$mdFile  = "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.md"
$zhXlf   = "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.zh-cn.xlf"
$deXlf   = "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.de-de.xlf"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/e74a80ce41411e17364d9d6491330c01720248cd/e2e/4f872ba6-5f5e-4fed-ae4b-08952f4241ec.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/edc167c13b7933011032f418179e713bfda96518/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dffc0685b23fede7a4f0130c09db4b2cde1f604d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.de-de.xlf"

# --- Overview sheet: refresh the per-locale status column text -----------
$ov.Range("B2").Value = $statusHandedBack
$ov.Range("C2").Value = $statusHandedBack
$ov.Range("B3").Value = $statusHandedBack
$ov.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet -----------------------------------------------------------
$zh.Range("B2").Value = $statusHandedBack
$zh.Hyperlinks.Add($zh.Range("E2"), $mdUrl, "", "", $mdFile)
$zh.Hyperlinks.Add($zh.Range("F2"), $zhXlfUrl, "", "", $zhXlf)
$zh.Range("G2").Value = "2016-01-25 04:02:19"
$zh.Range("H2").Value = "Include"

$zh.Range("B3").Value = $statusHandedBack
$zh.Hyperlinks.Add($zh.Range("E3"), $mdUrl, "", "", $mdFile)
$zh.Hyperlinks.Add($zh.Range("F3"), $zhXlfUrl, "", "", $zhXlf)
$zh.Range("G3").Value = "2016-01-25 04:02:19"
$zh.Range("H3").Value = "Include"

$zh.Range("D4").Value = "0001-01-01 00:00:00"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Ignored"

# --- de-de sheet -----------------------------------------------------------
$de.Range("B2").Value = $statusHandedBack
$de.Hyperlinks.Add($de.Range("E2"), $mdUrl, "", "", $mdFile)
$de.Hyperlinks.Add($de.Range("F2"), $deXlfUrl, "", "", $deXlf)
$de.Range("G2").Value = "2016-01-25 04:02:36"
$de.Range("H2").Value = "Include"

$de.Range("B3").Value = $statusHandedBack
$de.Hyperlinks.Add($de.Range("E3"), $mdUrl, "", "", $mdFile)
$de.Hyperlinks.Add($de.Range("F3"), $deXlfUrl, "", "", $deXlf)
$de.Range("G3").Value = "2016-01-25 04:02:36"
$de.Range("H3").Value = "Include"

$de.Range("D4").Value = "0001-01-01 00:00:00"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Ignored"

Write-Host "Handback report regenerated."
